# Updates to latest 4.0
# - "About" sheet: A4 ("None") loses its bold-ish style override (s="6" -> default)
# - "DPbES" sheet becomes the active/selected sheet, with a new selection
# - "DPbES" data: dispatch-priority flags flip for hard coal / hydro / biomass / petroleum

$wb = $excel.ActiveWorkbook

$about = $wb.Worksheets.Item("About")
$dpbes = $wb.Worksheets.Item("DPbES")

# --- About sheet: remove the extra font-applied style on A4 ("None") ---
$about.Range("A4").Font.Bold = $false

# --- DPbES sheet: flip the dispatch-priority flag rows ---
# Row 2  = hard coal   : 0 -> 1
$dpbes.Range("B2:AE2").Value = 1
# Row 6  = hydro       : 1 -> 0
$dpbes.Range("B6:AE6").Value = 0
# Row 10 = biomass     : 1 -> 0
$dpbes.Range("B10:AE10").Value = 0
# Row 12 = petroleum   : 0 -> 1
$dpbes.Range("B12:AE12").Value = 1

# --- Make DPbES the active sheet with the new selection/view state ---
$dpbes.Activate()
$dpbes.Range("B10:AE10").Select()
